# Refresh crypto price/volume snapshot (GitHub Actions scheduled update).
# Most D (Price) / E (Volume 1h) cells get new values; a few coin pairs
# (rows 36/37, 41/42, 50/51) swapped rank order with their own new figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.704.36"
$ws.Range("E2").Value = "  -4.36%  "
$ws.Range("D3").Value = "2.978.16"
$ws.Range("E3").Value = "  -5.23%  "
$ws.Range("D5").Value = "'540.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.65%  "
$ws.Range("D6").Value = "'151.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.77%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.76%  "
$ws.Range("D9").Value = "2.990.66"
$ws.Range("E9").Value = "  -5.42%  "
$ws.Range("E10").Value = "  -4.11%  "
$ws.Range("E11").Value = "  -7.34%  "
$ws.Range("D12").Value = "'0.368"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.27%  "
$ws.Range("D13").Value = "3.503.25"
$ws.Range("E13").Value = "  -5.37%  "
$ws.Range("D15").Value = "61.757.53"
$ws.Range("E15").Value = "  -4.33%  "
$ws.Range("D16").Value = "'23.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.48%  "
$ws.Range("D17").Value = "2.979.45"
$ws.Range("E17").Value = "  -5.54%  "
$ws.Range("E18").Value = "  -5.83%  "
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").Value = "'12.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.86%  "
$ws.Range("D21").Value = "'381.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.36%  "
$ws.Range("E22").Value = "  -5.24%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "'5.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.60%  "
$ws.Range("D25").Value = "'65.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.47%  "
$ws.Range("D26").Value = "'0.471"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("D27").Value = "3.104.31"
$ws.Range("E27").Value = "  -5.25%  "
$ws.Range("D28").Value = "'0.189"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.51%  "
$ws.Range("D29").Value = "'0.990"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "0.0₃0936"
$ws.Range("E30").Value = "  -8.51%  "
$ws.Range("D31").Value = "'8.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  -4.68%  "
$ws.Range("D34").Value = "'20.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.74%  "
$ws.Range("D35").Value = "'160.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.90%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "'5.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.01%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.38%  "
$ws.Range("E38").Value = "  -5.26%  "
$ws.Range("E39").Value = "  -6.56%  "
$ws.Range("D40").Value = "'1.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.63%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.424.13"
$ws.Range("E41").Value = "  -8.48%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'37.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("D43").Value = "'3.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.10%  "
$ws.Range("D44").Value = "'21.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.65%  "
$ws.Range("D45").Value = "'0.671"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.06%  "
$ws.Range("D47").Value = "'5.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.78%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").Value = "'0.0244"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.42%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'19.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.01%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "'0.0952"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.52%  "
